$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-09"

# Update the header label in I1 (year-to-date column header)
$ws.Range("I1").Value = "2022 (through 11-09)"

# Update November's 2022 figure (row 12 corresponds to November)
$ws.Range("I12").Value = 25

# Update the Total row's 2022 figure (row 14)
$ws.Range("I14").Value = 1424
